$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 updates
$ws.Range("D12").Value = 44489
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 8500
$ws.Range("P12").Value = 340

# Row 13 updates
$ws.Range("D13").Value = 44165
$ws.Range("J13").Value = 38
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 8500
$ws.Range("M13").Value = 8263
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 331
